$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "214.69") need an explicit
# text number format so Excel keeps them as text instead of converting to a number,
# matching how the source data is stored (inline/shared strings) in the workbook.
$textCells = @('D5', 'D6', 'D8', 'D10', 'D11', 'D14', 'D15', 'D18', 'D19', 'D24', 'D25', 'D26', 'D27', 'D28', 'D31', 'D32', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D44', 'D47', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated coin price / volume values exactly as published in the new snapshot.
$ws.Range('D2').Value = '29.931.36'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = '1.632.79'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '214.69'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').Value = '0.517'
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('D8').Value = '28.66'
$ws.Range('E8').Value = '  -1.86%  '
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('D10').Value = '0.0609'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').Value = '0.0905'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '1.866.01'
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').Value = '1.636.47'
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = '0.562'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').Value = '9.23'
$ws.Range('E15').Value = '  +13.56%  '
$ws.Range('D16').Value = '29.941.49'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('D18').Value = '64.18'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = '240.89'
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('E22').Value = '  +2.11%  '
$ws.Range('E23').Value = '  +2.85%  '
$ws.Range('D24').Value = '2.17'
$ws.Range('E24').Value = '  +2.96%  '
$ws.Range('D25').Value = '158.11'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').Value = '15.47'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').Value = '0.109'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').Value = '6.58'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  +1.78%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '3.39'
$ws.Range('E31').Value = '  +4.29%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '1.10'
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '1.429.39'
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('D35').Value = '1.65'
$ws.Range('E35').Value = '  +4.96%  '
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').Value = '2.75'
$ws.Range('E37').Value = '  -3.70%  '
$ws.Range('D38').Value = '2.28'
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').Value = '0.0170'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = '75.64'
$ws.Range('E40').Value = '  +12.00%  '
$ws.Range('D41').Value = '0.552'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('E42').Value = '  +2.23%  '
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('D44').Value = '0.0491'
$ws.Range('E44').Value = '  -1.42%  '
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('D47').Value = '50.99'
$ws.Range('E47').Value = '  -8.51%  '
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('D49').Value = '1.773.02'
$ws.Range('E49').Value = '  +1.89%  '
$ws.Range('E50').Value = '  +12.65%  '
$ws.Range('D51').Value = '90.49'
$ws.Range('E51').Value = '  +4.38%  '
